$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# 1. Table structure: grow from 3 to 5 columns.
#    Insert a new column before col 1 ("Semaines" / week column)
#    Insert a new column before col 4 ("Objectif de la semaine" column)
# ---------------------------------------------------------------------
$t.Columns.Add($t.Columns.Item(1)) | Out-Null
$t.Columns.Add($t.Columns.Item(4)) | Out-Null

# New column order is now: [Semaines(1), Dates(2), Taches(3), Objectif(4), Notes(5)]

# ---------------------------------------------------------------------
# 2. Row structure: insert two new rows so the 5-row table becomes 7,
#    then drop one now-superfluous trailing empty row to land on 6.
# ---------------------------------------------------------------------
# insert a fresh row right after row 2 (becomes the future row 3)
$t.Rows.Add($t.Rows.Item(3)) | Out-Null
# insert a fresh row right after the (shifted) old row 3, now row 4 (becomes future row 5)
$t.Rows.Add($t.Rows.Item(5)) | Out-Null
# table now has 7 rows; delete the last (still fully empty) row
$t.Rows.Item(7).Delete()

# ---------------------------------------------------------------------
# 3. Column widths (dxa values from the grid are expressed in points
#    for the COM Width property, so divide by 20).
# ---------------------------------------------------------------------
$widths = @(1696, 1272, 2596, 1701, 1797)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Width = $widths[$c - 1] / 20.0
    }
}

# ---------------------------------------------------------------------
# 4. Row 1 - header row: add the two new header labels.
# ---------------------------------------------------------------------
$t.Cell(1, 1).Range.Text = "Semaines"
$t.Cell(1, 4).Range.Text = "Objectif de la semaine"

# ---------------------------------------------------------------------
# 5. Row 2 - becomes the "SEMAINE 1" banner row.
# ---------------------------------------------------------------------
$t.Cell(2, 1).Range.Text = "SEMAINE 1 : apprentissage"
$t.Cell(2, 2).Range.Text = ""
$t.Cell(2, 3).Range.Text = ""
$t.Cell(2, 4).Range.Text = "se familiariser avec les outils et explorer les diff" + [char]0xE9 + "rentes fonctionnalit" + [char]0xE9 + "s offertes avec le robot."
# Cell(2,5) already holds "Debut du travail de Bachelor" - keep it as-is.

# bookmark _GoBack now lives alone in the (now empty) Cell(2,2)
$d.Bookmarks.Add("_GoBack", $t.Cell(2, 2).Range) | Out-Null

# shade the whole SEMAINE 1 row
for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell(2, $c)
    $cell.Shading.Texture = 0
    $cell.Shading.ForegroundPatternColor = -16777216
    $cell.Shading.BackgroundPatternColor = 0x0000C0
}

# ---------------------------------------------------------------------
# 6. Row 3 - new row holding the former "04-10-2016" entry.
# ---------------------------------------------------------------------
$t.Cell(3, 2).Range.Text = "04-10-2016"
$t.Cell(3, 3).Range.Text = "Installation de choregraphe et familiarisation avec l'outil"

# ---------------------------------------------------------------------
# 7. Row 4 - the old "06-10-2016" row; drop its bookmark (it now lives
#    on row 2) by rewriting the cell text without it. Content is
#    otherwise unchanged.
# ---------------------------------------------------------------------
$t.Cell(4, 3).Range.Text = "Exploration des fonctionnalit" + [char]0xE9 + "s du robot avec le simulateur(choregraphe) au moyens de s" + [char]0xE9 + "ries de tests"

# ---------------------------------------------------------------------
# 8. Row 5 - brand-new row for the "07-10-2016" entry.
# ---------------------------------------------------------------------
$t.Cell(5, 2).Range.Text = "07-10-2016"
$t.Cell(5, 3).Range.Text = "-Tests des fonctionnalit" + [char]0xE9 + "s avec Nao et Chroregraphe.`r"
$t.Cell(5, 5).Range.Text = "`r"

# ---------------------------------------------------------------------
# 9. Row 6 - becomes the "SEMAINE 2" banner row.
# ---------------------------------------------------------------------
$t.Cell(6, 1).Range.Text = "SEMAINE 2 : apprentissage 2"
$t.Cell(6, 4).Range.Text = "-Ecrire et Tester quelques sc" + [char]0xE9 + "narios sur le th" + [char]0xE8 + "me " + [char]0xAB + " Les portes ouvertes HEIG-VD " + [char]0xBB

for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell(6, $c)
    $cell.Shading.Texture = 0
    $cell.Shading.ForegroundPatternColor = -16777216
    $cell.Shading.BackgroundPatternColor = 0x0000C0
}
